# Format photo ID and select file .xlsx
#
# - H2/I2 header cells lose their (unused) "applyFill" variant style and
#   match the plain header style used elsewhere (same as A2, etc).
# - The date columns (H3:H8) keep their date format, just re-expressed
#   once the unused style above is retired.
# - Three new employees (rows 9-11) are filled in with name/role/photo-id
#   date/employee-number data; row 9's "B" column (second given name) has
#   no value for Mariana, so it is fully cleared (no cell record at all).
# - The active selection moves to I12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-style the H2 / I2 header cells to match the plain header style ---
$ws.Range("A2").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H2").Value = "Vigencia"
$ws.Range("I2").Value = "Numero de empleado"

# --- Row 9: Mariana Pacheco Diaz, Tecnico Docente ---
$ws.Range("A9").Value = "Mariana"
$ws.Range("B9").Clear()
$ws.Range("C9").Value = "Pacheco"
$ws.Range("D9").Value = "Diaz"
$ws.Range("E9").Value = "Técnico Docente"
$ws.Range("H3").Copy()
$ws.Range("H9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H9").Value = 45345
$ws.Range("I9").Value = 12345

# --- Row 10: Arturo Antonio Otal Morales, Tecnico Docente ---
$ws.Range("A10").Value = "Arturo"
$ws.Range("B10").Value = "Antonio"
$ws.Range("C10").Value = "Otal"
$ws.Range("D10").Value = "Morales"
$ws.Range("E10").Value = "Técnico Docente"
$ws.Range("H3").Copy()
$ws.Range("H10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H10").Value = 45345
$ws.Range("I10").Value = 12345

# --- Row 11: Karla Guadalupe Nava Rodríguez, Gestoria ---
$ws.Range("A11").Value = "Karla"
$ws.Range("B11").Value = "Guadalupe"
$ws.Range("C11").Value = "Nava"
$ws.Range("D11").Value = "Rodríguez"
$ws.Range("E11").Value = "Gestoría"
$ws.Range("H3").Copy()
$ws.Range("H11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H11").Value = 45410
$ws.Range("I11").Value = 12345

# --- Move the active selection to I12, matching the committed file ---
$ws.Range("I12").Select() | Out-Null
